# Update the two-digit division worksheet numbers.
# Each Find/Replace targets the first remaining occurrence (wdReplaceOne = 1)
# in document order, so duplicate "old" values (e.g. "42÷3=" which appears
# twice, with two different replacements) resolve correctly.

$d = $word.ActiveDocument

$replacements = @(
    @("97÷5=", "27÷8="),
    @("70÷7=", "49÷9="),
    @("63÷9=", "78÷6="),
    @("74÷6=", "52÷3="),
    @("67÷2=", "97÷7="),
    @("45÷7=", "21÷7="),
    @("13÷3=", "99÷6="),
    @("28÷3=", "38÷4="),
    @("76÷2=", "98÷7="),
    @("93÷5=", "24÷9="),
    @("42÷3=", "50÷8="),
    @("79÷2=", "97÷5="),
    @("69÷8=", "66÷2="),
    @("68÷2=", "51÷5="),
    @("73÷9=", "45÷2="),
    @("31÷6=", "80÷5="),
    @("33÷3=", "96÷4="),
    @("69÷4=", "32÷4="),
    @("42÷3=", "21÷3="),
    @("91÷2=", "34÷4="),
    @("62÷6=", "24÷3="),
    @("94÷8=", "29÷8="),
    @("64÷9=", "84÷9="),
    @("18÷5=", "99÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 1)
}
